$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B23").Value = 6333
$ws.Range("C23").Value = 1003
$ws.Range("D23").Value = 5908082
$ws.Range("E23").Value = 932.90415285015
$ws.Range("F23").Value = 8.665065202470835
$ws.Range("G23").Value = 4.370447450572312
$ws.Range("H23").Value = 26.59507473141238
